$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 63
$ws.Range("A3").Value = 64
$ws.Range("A4").Value = 65
$ws.Range("A5").Value = 66
